$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 32-39 with shifted weekly price data ---
# Row 32
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 55
$ws.Range("K32").Value = 2600
$ws.Range("L32").Value = 2600
$ws.Range("M32").Value = 2600
$ws.Range("P32").Value = 2600

# Row 33
$ws.Range("D33").Value = 44159
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H33").Value = "Verde"
$ws.Range("I33").Value = "Banquete"
$ws.Range("J33").Value = 180
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 10000
$ws.Range("N33").Value = "$/bandeja 10 kilos"
$ws.Range("P33").Value = 1000
$ws.Range("Q33").Value = 10

# Row 34
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 260
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = 9000
$ws.Range("P34").Value = 900

# Row 35
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 320
$ws.Range("K35").Value = 7000
$ws.Range("L35").Value = 7000
$ws.Range("M35").Value = 7000
$ws.Range("P35").Value = 700

# Row 36
$ws.Range("D36").Value = 44476
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Banquete"
$ws.Range("J36").Value = 180
$ws.Range("K36").Value = 1400
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = 1456
$ws.Range("N36").Value = "$/kilo"
$ws.Range("O36").Value = "Región Metropolitana"
$ws.Range("P36").Value = 1456
$ws.Range("Q36").Value = 1

# Row 37
$ws.Range("D37").Value = 44476
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 270
$ws.Range("K37").Value = 1100
$ws.Range("L37").Value = 1200
$ws.Range("M37").Value = 1144
$ws.Range("P37").Value = 1144

# Row 38
$ws.Range("D38").Value = 44476
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value = 370
$ws.Range("K38").Value = 900
$ws.Range("L38").Value = 1000
$ws.Range("M38").Value = 959
$ws.Range("P38").Value = 959

# Row 39
$ws.Range("I39").Value = "Banquete"
$ws.Range("J39").Value = 90
$ws.Range("K39").Value = 1600
$ws.Range("L39").Value = 1600
$ws.Range("M39").Value = 1600
$ws.Range("P39").Value = 1600

# --- Append two new rows (40, 41) for the latest week ---
# Row 40
$ws.Range("A40").Value = 12
$ws.Range("B40").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 44468
$ws.Range("D40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = 300000000
$ws.Range("G40").Value = "Espárragos"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 80
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = 1500
$ws.Range("N40").Value = "$/kilo"
$ws.Range("O40").Value = "Región Metropolitana"
$ws.Range("P40").Value = 1500
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = "Hortaliza"

# Row 41
$ws.Range("A41").Value = 12
$ws.Range("B41").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C41").Value = "Metropolitana"
$ws.Range("D41").Value = 44468
$ws.Range("D41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = 300000000
$ws.Range("G41").Value = "Espárragos"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Segunda"
$ws.Range("J41").Value = 75
$ws.Range("K41").Value = 1300
$ws.Range("L41").Value = 1300
$ws.Range("M41").Value = 1300
$ws.Range("N41").Value = "$/kilo"
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 1300
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"

